$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3348.5
$ws.Range("J17").Value = 3401.578
$ws.Range("L17").Value = 10204.734
$ws.Range("N17").Value = -10540.734

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 23810.682
$ws.Range("I64").Value = 93227.27
$ws.Range("J64").Value = 2600.0557
$ws.Range("K64").Value = 93227.27
$ws.Range("L64").Value = 2600.0557
$ws.Range("M64").Value = -92979.27
$ws.Range("N64").Value = -3096.0557

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 23810.682
$ws.Range("I67").Value = 93227.27
$ws.Range("J67").Value = 2600.0557
$ws.Range("K67").Value = 93227.27
$ws.Range("L67").Value = 2600.0557
$ws.Range("M67").Value = -92369.27
$ws.Range("N67").Value = -4316.0557

# ALC row 123
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 413.10345
$ws.Range("I125").Value = 778.2857
$ws.Range("J125").Value = 296.9091
$ws.Range("K125").Value = 7004.571300000001
$ws.Range("L125").Value = 2672.1819
$ws.Range("M125").Value = -4544.571300000001
$ws.Range("N125").Value = -7592.1819

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2017.4342
$ws.Range("I138").Value = 2475.0588
$ws.Range("J138").Value = 1885.5763
$ws.Range("K138").Value = 7425.176399999999
$ws.Range("L138").Value = 5656.7289
$ws.Range("M138").Value = -2285.176399999999
$ws.Range("N138").Value = -15936.7289

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 451.66666
$ws.Range("I5").Value = 385.7143
$ws.Range("J5").Value = 544
$ws.Range("K5").Value = 385.7143
$ws.Range("L5").Value = 544
$ws.Range("M5").Value = -273.7143
$ws.Range("N5").Value = -768

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2722.926
$ws.Range("I63").Value = 2409.9546
$ws.Range("J63").Value = 4100
$ws.Range("K63").Value = 2409.9546
$ws.Range("L63").Value = 4100
$ws.Range("M63").Value = -1723.9546
$ws.Range("N63").Value = -5472

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2722.926
$ws.Range("I66").Value = 2409.9546
$ws.Range("J66").Value = 4100
$ws.Range("K66").Value = 12049.773
$ws.Range("L66").Value = 20500
$ws.Range("M66").Value = -8617.773000000001
$ws.Range("N66").Value = -27364

# ARM row 107
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 38493.5
$ws.Range("J107").Value = 38493.5
$ws.Range("L107").Value = 38493.5
$ws.Range("N107").Value = -46173.5

# ARM row 109
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 44988
$ws.Range("J109").Value = 44988
$ws.Range("L109").Value = 44988
$ws.Range("N109").Value = -47762

# ARM row 117
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H117").Value = 40230.57
$ws.Range("J117").Value = 40230.57
$ws.Range("L117").Value = 40230.57
$ws.Range("N117").Value = -49408.57

# ARM row 119
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 36403.5
$ws.Range("J119").Value = 36403.5
$ws.Range("L119").Value = 36403.5
$ws.Range("N119").Value = -46079.5

# ARM row 137
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 31627.25
$ws.Range("J137").Value = 41933.332
$ws.Range("L137").Value = 41933.332
$ws.Range("N137").Value = -52133.332

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 451.66666
$ws.Range("I4").Value = 385.7143
$ws.Range("J4").Value = 544
$ws.Range("K4").Value = 385.7143
$ws.Range("L4").Value = 544
$ws.Range("M4").Value = -270.7143
$ws.Range("N4").Value = -774

# BSM row 122
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 40566
$ws.Range("J122").Value = 40566
$ws.Range("L122").Value = 40566
$ws.Range("N122").Value = -50366

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2835.5881
$ws.Range("I62").Value = 2728.5
$ws.Range("J62").Value = 3335.3333
$ws.Range("K62").Value = 2728.5
$ws.Range("L62").Value = 3335.3333
$ws.Range("M62").Value = -2104.5
$ws.Range("N62").Value = -4583.3333

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2835.5881
$ws.Range("I65").Value = 2728.5
$ws.Range("J65").Value = 3335.3333
$ws.Range("K65").Value = 13642.5
$ws.Range("L65").Value = 16676.6665
$ws.Range("M65").Value = -10522.5
$ws.Range("N65").Value = -22916.6665

# CRP row 80
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 30081.9
$ws.Range("J80").Value = 30081.9
$ws.Range("L80").Value = 30081.9
$ws.Range("N80").Value = -32327.9

# CRP row 83
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H83").Value = 30081.9
$ws.Range("J83").Value = 30081.9
$ws.Range("L83").Value = 90245.70000000001
$ws.Range("N83").Value = -101477.7

# CRP row 104
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H104").Value = 29997
$ws.Range("J104").Value = 29997
$ws.Range("L104").Value = 29997
$ws.Range("N104").Value = -35239

# CRP row 109
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 27657.4
$ws.Range("J109").Value = 27657.4
$ws.Range("L109").Value = 27657.4
$ws.Range("N109").Value = -29737.4

# CRP row 111
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H111").Value = 48398.668
$ws.Range("J111").Value = 48398.668
$ws.Range("L111").Value = 48398.668
$ws.Range("N111").Value = -56578.668

# CRP row 112
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 32017.8
$ws.Range("J112").Value = 32017.8
$ws.Range("L112").Value = 32017.8
$ws.Range("N112").Value = -34971.8

# CRP row 115
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H115").Value = 28919.334
$ws.Range("J115").Value = 28919.334
$ws.Range("L115").Value = 28919.334
$ws.Range("N115").Value = -31269.334

# CRP row 116
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 40918.75
$ws.Range("J116").Value = 40918.75
$ws.Range("L116").Value = 40918.75
$ws.Range("N116").Value = -50096.75

# CRP row 118
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H118").Value = 41420
$ws.Range("J118").Value = 41420
$ws.Range("L118").Value = 41420
$ws.Range("N118").Value = -44734

# CRP row 119
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H119").Value = 42628.25
$ws.Range("J119").Value = 42628.25
$ws.Range("L119").Value = 42628.25
$ws.Range("N119").Value = -52304.25

# CRP row 120
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H120").Value = 29115.188
$ws.Range("J120").Value = 29115.188
$ws.Range("L120").Value = 29115.188
$ws.Range("N120").Value = -36373.18799999999

# CRP row 133
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 25209.092
$ws.Range("J133").Value = 25209.092
$ws.Range("L133").Value = 25209.092
$ws.Range("N133").Value = -30269.092

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3247.342
$ws.Range("J113").Value = 657
$ws.Range("L113").Value = 1971
$ws.Range("N113").Value = -6311

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 40986.5
$ws.Range("I140").Value = 92521
$ws.Range("J140").Value = 3194.5334
$ws.Range("K140").Value = 277563
$ws.Range("L140").Value = 9583.600199999999
$ws.Range("M140").Value = -272383
$ws.Range("N140").Value = -19943.6002

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 150785.6
$ws.Range("I80").Value = 266832.3
$ws.Range("J80").Value = 3793.0667
$ws.Range("K80").Value = 266832.3
$ws.Range("L80").Value = 3793.0667
$ws.Range("M80").Value = -265834.3
$ws.Range("N80").Value = -5789.066699999999

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 150785.6
$ws.Range("I83").Value = 266832.3
$ws.Range("J83").Value = 3793.0667
$ws.Range("K83").Value = 1334161.5
$ws.Range("L83").Value = 18965.3335
$ws.Range("M83").Value = -1329169.5
$ws.Range("N83").Value = -28949.3335

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1677.3889
$ws.Range("I113").Value = 1731.1111
$ws.Range("K113").Value = 1731.1111
$ws.Range("M113").Value = 438.8888999999999
$ws.Range("N113").ClearContents()

# GSM row 116
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H116").Value = 38996
$ws.Range("J116").Value = 38996
$ws.Range("L116").Value = 38996
$ws.Range("N116").Value = -48174

# GSM row 118
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 38086.668
$ws.Range("J118").Value = 38086.668
$ws.Range("L118").Value = 38086.668
$ws.Range("N118").Value = -41400.668

# GSM row 130
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 45850
$ws.Range("J130").Value = 45850
$ws.Range("L130").Value = 45850
$ws.Range("N130").Value = -55890

# LTW row 9
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 7269.3335
$ws.Range("I9").Value = 404
$ws.Range("J9").Value = 21000
$ws.Range("K9").Value = 404
$ws.Range("L9").Value = 21000
$ws.Range("M9").Value = -180
$ws.Range("N9").Value = -21448

# LTW row 20
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 269.73077
$ws.Range("I55").Value = 185.61539
$ws.Range("J55").Value = 353.84616
$ws.Range("K55").Value = 185.61539
$ws.Range("L55").Value = 353.84616
$ws.Range("M55").Value = -12.61538999999999
$ws.Range("N55").Value = -699.8461600000001

# LTW row 108
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 48614
$ws.Range("J108").Value = 48614
$ws.Range("L108").Value = 48614
$ws.Range("N108").Value = -56294

# LTW row 119
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H119").Value = 47408
$ws.Range("J119").Value = 47408
$ws.Range("L119").Value = 47408
$ws.Range("N119").Value = -57084

# LTW row 137
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 25006.666
$ws.Range("J137").Value = 25006.666
$ws.Range("L137").Value = 25006.666
$ws.Range("N137").Value = -35206.666

# LTW row 139
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 44721.355
$ws.Range("J139").Value = 44721.355
$ws.Range("L139").Value = 44721.355
$ws.Range("N139").Value = -55001.355

# WVR row 16
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 45925.75
$ws.Range("J16").Value = 45925.75
$ws.Range("L16").Value = 45925.75
$ws.Range("N16").Value = -46509.75

# WVR row 119
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 48694
$ws.Range("J119").Value = 48694
$ws.Range("L119").Value = 48694
$ws.Range("N119").Value = -58370

# WVR row 120
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H120").Value = 45416
$ws.Range("J120").Value = 45416
$ws.Range("L120").Value = 45416
$ws.Range("N120").Value = -55092

# WVR row 121
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 38447.2
$ws.Range("J121").Value = 38447.2
$ws.Range("L121").Value = 38447.2
$ws.Range("N121").Value = -41941.2

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 14286714
$ws.Range("I122").Value = 14286714
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 42860142
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -42857692
$ws.Range("N122").ClearContents()

# WVR row 139
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 20801.154
$ws.Range("J139").Value = 20801.154
$ws.Range("L139").Value = 20801.154
$ws.Range("N139").Value = -31081.154
